$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Endpoints")

# Insert a new row at the top of the "Endpoints" sheet, shifting every
# existing row (and the merged A-column cells / row references that go
# with them) down by one.
$ws.Rows.Item(1).Insert()

# Put the new title text in A1, then merge A1:F1 into a single banner cell.
$ws.Range("A1").Value = "Defining the Endpoints"
$titleRange = $ws.Range("A1:F1")
$titleRange.Merge()

# Style the merged title: centered, bold 14pt Roboto, with a thin bottom
# border, and a taller row to fit it.
$titleRange.HorizontalAlignment = -4108  # xlCenter
$titleRange.Font.Bold = $true
$titleRange.Font.Size = 14
$titleRange.Font.Name = "Roboto"
$titleRange.Borders.Item(9).LineStyle = 1  # xlEdgeBottom
$titleRange.Borders.Item(9).Weight = 2     # xlThin
$ws.Rows.Item(1).RowHeight = 18

# Make "Endpoints" the active sheet/tab, with the new title selected.
$ws.Activate()
$titleRange.Select()

$wb.Save()
